$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "style"
$ws.Range("A3").Value = "c-sigma"
